$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "Rooms"

# Update header labels
$ws.Range("A1").Value = "room_num"
$ws.Range("B1").Value = "room_type"

# Fill in the Capacity column for the previously-empty rows (2-8)
$ws.Range("E2:E8").Value = 10

# Update the selection shown when the workbook is opened
$ws.Range("E2:E8").Select()
